$wb = $excel.ActiveWorkbook

# ---- Sheet ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H21").Value = 37000
$ws.Range("I21").Value = 37000
$ws.Range("K21").Value = 37000
$ws.Range("M21").Value = -36532
$ws.Range("H23").Value = 37000
$ws.Range("I23").Value = 37000
$ws.Range("K23").Value = 37000
$ws.Range("M23").Value = -36766
$ws.Range("H33").Value = 285.57144
$ws.Range("J33").Value = 626
$ws.Range("L33").Value = 626
$ws.Range("N33").Value = -1084
$ws.Range("H42").Value = 2512.5
$ws.Range("I42").Value = 75
$ws.Range("K42").Value = 225
$ws.Range("M42").Value = 5
$ws.Range("H59").Value = 118
$ws.Range("J59").Value = 118
$ws.Range("L59").Value = 354
$ws.Range("N59").Value = -1468
$ws.Range("H132").Value = 4940.8335
$ws.Range("I132").Value = 4909.8
$ws.Range("J132").Value = 5096
$ws.Range("K132").Value = 14729.4
$ws.Range("L132").Value = 15288
$ws.Range("M132").Value = -12199.4
$ws.Range("N132").Value = -20348
$ws.Range("H135").Value = 2027.9
$ws.Range("I135").Value = 2042.1111
$ws.Range("K135").Value = 18378.9999
$ws.Range("M135").Value = -15843.9999
$ws.Range("H141").Value = 33795
$ws.Range("I141").Value = 33795
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 101385
$ws.Range("L141").Value = 0
$ws.Range("M141").ClearContents()
$ws.Range("N141").Value = -96205

# ---- Sheet ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H19").Value = 504
$ws.Range("I19").Value = 504
$ws.Range("J19").Value = 0
$ws.Range("K19").Value = 504
$ws.Range("L19").ClearContents()
$ws.Range("N19").Value = 0
$ws.Range("M19").Value = -275
$ws.Range("H45").Value = 2315.5557
$ws.Range("I45").Value = 2503.5
$ws.Range("J45").Value = 1939.6666
$ws.Range("K45").Value = 2503.5
$ws.Range("L45").Value = 1939.6666
$ws.Range("M45").Value = -2126.5
$ws.Range("N45").Value = -2693.6666
$ws.Range("H102").Value = 7810
$ws.Range("I102").Value = 7810
$ws.Range("K102").Value = 7810
$ws.Range("M102").Value = -6188
$ws.Range("H122").Value = 4475
$ws.Range("I122").Value = 4475
$ws.Range("K122").Value = 13425
$ws.Range("M122").Value = -10975
$ws.Range("H132").Value = 0
$ws.Range("I132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("M132").ClearContents()

# ---- Sheet BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3429.25
$ws.Range("I20").Value = 2589
$ws.Range("K20").Value = 2589
$ws.Range("M20").Value = -2342
$ws.Range("H134").Value = 2187.4285
$ws.Range("I134").Value = 2235
$ws.Range("J134").Value = 1902
$ws.Range("K134").Value = 6705
$ws.Range("L134").Value = 5706
$ws.Range("M134").Value = -4170
$ws.Range("N134").Value = -10776

# ---- Sheet CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H32").Value = 337999.66
$ws.Range("I32").Value = 337999.66
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 337999.66
$ws.Range("L32").ClearContents()
$ws.Range("N32").Value = 0
$ws.Range("M32").Value = -337683.66
$ws.Range("H121").Value = 0
$ws.Range("J121").Value = 0
$ws.Range("L121").ClearContents()
$ws.Range("N121").Value = 0
$ws.Range("H132").Value = 3327.111
$ws.Range("I132").Value = 2688.8
$ws.Range("K132").Value = 8066.400000000001
$ws.Range("M132").Value = -5536.400000000001

# ---- Sheet CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 2000
$ws.Range("I3").Value = 0
$ws.Range("K3").Value = 0
$ws.Range("M3").ClearContents()
$ws.Range("H4").Value = 441.9091
$ws.Range("I4").Value = 94.53846
$ws.Range("J4").Value = 943.6667
$ws.Range("K4").Value = 283.61538
$ws.Range("L4").Value = 2831.0001
$ws.Range("M4").Value = -171.61538
$ws.Range("N4").Value = -3055.0001
$ws.Range("H121").Value = 1309.9524
$ws.Range("I121").Value = 1030.75
$ws.Range("K121").Value = 3092.25
$ws.Range("M121").Value = -1782.25

# ---- Sheet GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 19.333334
$ws.Range("I2").Value = 16.5
$ws.Range("K2").Value = 16.5
$ws.Range("M2").Value = 96.5
$ws.Range("H14").Value = 1050045
$ws.Range("I14").Value = 1312543.8
$ws.Range("J14").Value = 50
$ws.Range("K14").Value = 1312543.8
$ws.Range("L14").Value = 50
$ws.Range("M14").Value = -1312375.8
$ws.Range("N14").Value = -386
$ws.Range("H20").Value = 3338666.8
$ws.Range("I20").Value = 5000500
$ws.Range("J20").Value = 15000
$ws.Range("K20").Value = 5000500
$ws.Range("L20").Value = 15000
$ws.Range("M20").Value = -5000255
$ws.Range("N20").Value = -15490
$ws.Range("H57").Value = 16030.5
$ws.Range("I57").Value = 0
$ws.Range("J57").Value = 16030.5
$ws.Range("K57").Value = 0
$ws.Range("L57").ClearContents()
$ws.Range("M57").Value = 16030.5
$ws.Range("N57").Value = -17670.5
$ws.Range("H80").Value = 9281
$ws.Range("I80").Value = 2567.25
$ws.Range("J80").Value = 18232.666
$ws.Range("K80").Value = 2567.25
$ws.Range("L80").Value = 18232.666
$ws.Range("M80").Value = -1569.25
$ws.Range("N80").Value = -20228.666
$ws.Range("H83").Value = 9281
$ws.Range("I83").Value = 2567.25
$ws.Range("J83").Value = 18232.666
$ws.Range("K83").Value = 12836.25
$ws.Range("L83").Value = 91163.33
$ws.Range("M83").Value = -7844.25
$ws.Range("N83").Value = -101147.33
$ws.Range("H113").Value = 1350
$ws.Range("I113").Value = 1200
$ws.Range("K113").Value = 1200
$ws.Range("M113").Value = 970

# ---- Sheet LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1960.6
$ws.Range("I7").Value = 1934.3334
$ws.Range("K7").Value = 1934.3334
$ws.Range("M7").Value = -1822.3334
$ws.Range("H16").Value = 3249.8823
$ws.Range("J16").Value = 3800
$ws.Range("L16").Value = 3800
$ws.Range("N16").Value = -4140
$ws.Range("H55").Value = 1036.75
$ws.Range("I55").Value = 1199.6666
$ws.Range("J55").Value = 548
$ws.Range("K55").Value = 1199.6666
$ws.Range("L55").Value = 548
$ws.Range("M55").Value = -1026.6666
$ws.Range("N55").Value = -894
$ws.Range("H68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("L68").ClearContents()
$ws.Range("N68").Value = 0
$ws.Range("H71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("L71").ClearContents()
$ws.Range("N71").Value = 0
$ws.Range("H122").Value = 7078.643
$ws.Range("I122").Value = 7078.643
$ws.Range("K122").Value = 21235.929
$ws.Range("M122").Value = -18785.929
$ws.Range("H126").Value = 1960.6
$ws.Range("I126").Value = 1934.3334
$ws.Range("K126").Value = 5803.0002
$ws.Range("M126").Value = -3333.0002
$ws.Range("H136").Value = 0
$ws.Range("I136").Value = 0
$ws.Range("K136").Value = 0
$ws.Range("M136").ClearContents()

# ---- Sheet WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 8667.666999999999
$ws.Range("J2").Value = 3001.5
$ws.Range("L2").Value = 3001.5
$ws.Range("N2").Value = -3225.5
$ws.Range("H126").Value = 1511.2174
$ws.Range("I126").Value = 1359.7
$ws.Range("J126").Value = 2521.3333
$ws.Range("K126").Value = 4079.1
$ws.Range("L126").Value = 7563.999899999999
$ws.Range("M126").Value = -1609.1
